$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# Rename the existing sheet
$ws1.Name = "Chau Hoang"

# Add a new sheet right after it and give it the Vietnamese-accented name
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Châu Hoàng"

# Populate the new sheet with the same layout as sheet1 (Section / sub-items)
# but with a different set of numbers in column B
$ws2.Range("A1").Value = "Section"
$ws2.Range("B1").Value = "Tình yêu"
$ws2.Range("A2").Value = 1.1
$ws2.Range("B2").Value = 2
$ws2.Range("A3").Value = 1.2
$ws2.Range("B3").Value = 1
$ws2.Range("A4").Value = 2.1
$ws2.Range("B4").Value = 1
$ws2.Range("A5").Value = "Section"
$ws2.Range("B5").Value = "Tình bạn"
$ws2.Range("A6").Value = 3.1
$ws2.Range("B6").Value = 1
$ws2.Range("A7").Value = 3.2
$ws2.Range("B7").Value = 2
$ws2.Range("A8").Value = 4.1
$ws2.Range("B8").Value = 1
$ws2.Range("A9").Value = 2
$ws2.Range("B9").Value = 1

# Row heights: the two "Section" header rows are a touch taller than data rows
$ws2.Rows.Item(1).RowHeight = 12.65
$ws2.Rows.Item(2).RowHeight = 12.1
$ws2.Rows.Item(3).RowHeight = 12.1
$ws2.Rows.Item(4).RowHeight = 12.1
$ws2.Rows.Item(5).RowHeight = 12.65
$ws2.Rows.Item(6).RowHeight = 12.1
$ws2.Rows.Item(7).RowHeight = 12.1
$ws2.Rows.Item(8).RowHeight = 12.1
$ws2.Rows.Item(9).RowHeight = 12.1

# Sheet1 is no longer the tab-selected sheet and its selection resets to A1
$ws1.Range("A1").Select()

# The new sheet becomes active / tab-selected, with B9 selected
$ws2.Activate()
$ws2.Range("B9").Select()
